$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$url = "https://www.zinefilos.com/"
$desc = "Información de peliculas para usar data e imágenes"

$ws.Hyperlinks.Add($ws.Range("B8"), $url)
$ws.Range("B8").Style = "Hipervínculo"
$ws.Range("C8").Value = $desc

$ws.Range("C9").Select()
